$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: 'Bitcoin'
$ws.Range("D2").Value = "62.748.12"
$ws.Range("E2").Value = "  +6.56%  "

# Row 3: 'Ethereum'
$ws.Range("D3").Value = "3.468.17"
$ws.Range("E3").Value = "  +4.93%  "

# Row 4: 'TetherUSD'
$ws.Range("E4").Value = "  +0.16%  "

# Row 5: 'BNB'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "413.78"
$ws.Range("E5").Value = "  +3.26%  "

# Row 6: 'Solana'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "131.08"
$ws.Range("E6").Value = "  +19.18%  "

# Row 7: 'LidoStakedEther'
$ws.Range("D7").Value = "3.460.31"
$ws.Range("E7").Value = "  +4.79%  "

# Row 8: 'XRP'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.593"
$ws.Range("E8").Value = "  +0.82%  "

# Row 9: 'USDC'
$ws.Range("E9").Value = "  +0.13%  "

# Row 10: 'Cardano'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.690"
$ws.Range("E10").Value = "  +8.27%  "

# Row 11: 'Dogecoin'
$ws.Range("E11").Value = "  +29.42%  "

# Row 12: 'Avalanche'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "42.84"
$ws.Range("E12").Value = "  +7.45%  "

# Row 13: 'TRON'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.143"
$ws.Range("E13").Value = "  +0.35%  "

# Row 14: 'WrappedliquidstakedEther2.0'
$ws.Range("D14").Value = "4.018.17"
$ws.Range("E14").Value = "  +5.15%  "

# Row 15: 'Polkadot'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.75"
$ws.Range("E15").Value = "  +3.98%  "

# Row 16: 'Chainlink'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.15"
$ws.Range("E16").Value = "  +4.34%  "

# Row 17: 'WrappedEther'
$ws.Range("D17").Value = "3.458.90"
$ws.Range("E17").Value = "  +5.01%  "

# Row 18: 'WrappedBTC'
$ws.Range("D18").Value = "62.642.56"
$ws.Range("E18").Value = "  +7.02%  "

# Row 19: 'Polygon'
$ws.Range("E19").Value = "  +0.16%  "

# Row 20: 'Uniswap'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.87"
$ws.Range("E20").Value = "  -0.17%  "

# Row 21: 'ShibaInu'
$ws.Range("E21").Value = "  +25.33%  "

# Row 22: 'ImmutableX'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.38"
$ws.Range("E22").Value = "  +1.63%  "

# Row 23: 'Litecoin' -> 'InternetComputer(DFINITY)'
$ws.Range("B23").Value = "InternetComputer(DFINITY)"
$ws.Range("C23").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.36"
$ws.Range("E23").Value = "  +2.71%  "

# Row 24: 'BitcoinCash' -> 'Litecoin'
$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "82.31"
$ws.Range("E24").Value = "  +10.06%  "

# Row 25: 'InternetComputer(DFINITY)' -> 'BitcoinCash'
$ws.Range("B25").Value = "BitcoinCash"
$ws.Range("C25").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "316.01"
$ws.Range("E25").Value = "  +3.87%  "

# Row 26: 'PancakeSwap'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.19"
$ws.Range("E26").Value = "  +0.12%  "

# Row 27: 'EthereumClassic'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "30.72"
$ws.Range("E27").Value = "  +8.13%  "

# Row 28: 'Filecoin'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.15"
$ws.Range("E28").Value = "  +3.62%  "

# Row 29: 'RenderToken'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.69"
$ws.Range("E29").Value = "  +3.96%  "

# Row 30: 'LEO' -> 'Kaspa'
$ws.Range("B30").Value = "Kaspa"
$ws.Range("C30").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.178"
$ws.Range("E30").Value = "  +4.47%  "

# Row 31: 'Kaspa' -> 'LEO'
$ws.Range("B31").Value = "LEO"
$ws.Range("C31").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.39"
$ws.Range("E31").Value = "  -0.87%  "

# Row 32: 'Hedera'
$ws.Range("E32").Value = "  +4.80%  "

# Row 33: 'Toncoin' -> 'InjectiveProtocol'
$ws.Range("B33").Value = "InjectiveProtocol"
$ws.Range("C33").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "43.87"
$ws.Range("E33").Value = "  +9.54%  "

# Row 34: 'Cosmos' -> 'Toncoin'
$ws.Range("B34").Value = "Toncoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.64"
$ws.Range("E34").Value = "  +23.73%  "

# Row 35: 'InjectiveProtocol' -> 'Cosmos'
$ws.Range("B35").Value = "Cosmos"
$ws.Range("C35").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "11.80"
$ws.Range("E35").Value = "  +3.29%  "

# Row 36: 'Dai'
$ws.Range("E36").Value = "  -0.02%  "

# Row 37: 'VeChain'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0495"
$ws.Range("E37").Value = "  -6.33%  "

# Row 38: 'OKB'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "52.36"
$ws.Range("E38").Value = "  +0.91%  "

# Row 39: 'LidoDAOToken'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.58"
$ws.Range("E39").Value = "  +2.83%  "

# Row 40: 'FirstDigitalUSD'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.998"
$ws.Range("E40").Value = "  +0.05%  "

# Row 41: 'Stacks'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.04"
$ws.Range("E41").Value = "  -9.87%  "

# Row 42: 'ARBITRUM'
$ws.Range("E42").Value = "  +6.83%  "

# Row 43: 'Stellar'
$ws.Range("E43").Value = "  +2.77%  "

# Row 44: 'Monero'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "137.25"
$ws.Range("E44").Value = "  -0.31%  "

# Row 45: 'Celestia'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "17.35"
$ws.Range("E45").Value = "  +3.20%  "

# Row 46: 'TheGraph' -> 'NEARProtocol'
$ws.Range("B46").Value = "NEARProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.00"
$ws.Range("E46").Value = "  +1.39%  "

# Row 47: 'NEARProtocol' -> 'TheGraph'
$ws.Range("B47").Value = "TheGraph"
$ws.Range("C47").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.287"
$ws.Range("E47").Value = "  +2.86%  "

# Row 48: 'WEMIXToken'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.24"
$ws.Range("E48").Value = "  -3.20%  "

# Row 49: 'EnergySwap'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "22.19"
$ws.Range("E49").Value = "  -1.12%  "

# Row 50: 'Maker'
$ws.Range("D50").Value = "2.228.12"
$ws.Range("E50").Value = "  +2.87%  "

# Row 51: 'RocketPoolETH'
$ws.Range("D51").Value = "3.817.12"
$ws.Range("E51").Value = "  +5.18%  "

# Remove temporary text-number-format overrides so styling matches original (no explicit style on these cells)
$ws.Range("D5").ClearFormats()
$ws.Range("D6").ClearFormats()
$ws.Range("D8").ClearFormats()
$ws.Range("D10").ClearFormats()
$ws.Range("D12").ClearFormats()
$ws.Range("D13").ClearFormats()
$ws.Range("D15").ClearFormats()
$ws.Range("D16").ClearFormats()
$ws.Range("D20").ClearFormats()
$ws.Range("D22").ClearFormats()
$ws.Range("D23").ClearFormats()
$ws.Range("D24").ClearFormats()
$ws.Range("D25").ClearFormats()
$ws.Range("D26").ClearFormats()
$ws.Range("D27").ClearFormats()
$ws.Range("D28").ClearFormats()
$ws.Range("D29").ClearFormats()
$ws.Range("D30").ClearFormats()
$ws.Range("D31").ClearFormats()
$ws.Range("D33").ClearFormats()
$ws.Range("D34").ClearFormats()
$ws.Range("D35").ClearFormats()
$ws.Range("D37").ClearFormats()
$ws.Range("D38").ClearFormats()
$ws.Range("D39").ClearFormats()
$ws.Range("D40").ClearFormats()
$ws.Range("D41").ClearFormats()
$ws.Range("D44").ClearFormats()
$ws.Range("D45").ClearFormats()
$ws.Range("D46").ClearFormats()
$ws.Range("D47").ClearFormats()
$ws.Range("D48").ClearFormats()
$ws.Range("D49").ClearFormats()
